# Added routine to compress the DBs
# - Adds a new "sde" worksheet with user/password style data (mirrors the email sheet layout)
# - Adds a new "emailRecipient" worksheet listing recipient addresses
# - Makes "emailRecipient" the active/selected tab
# - Clears the old "A3" selection leftover on the "email" sheet, selecting A1:B2 instead

$wb = $excel.ActiveWorkbook

# --- New sheet: sde ---
$sdeSheet = $wb.Worksheets.Add()
$sdeSheet.Name = "sde"
$sdeSheet.Range("A1").Value = "user"
$sdeSheet.Range("B1").Value = "password"
$sdeSheet.Range("A2").Value = "sde"
$sdeSheet.Range("B2").Value = "*&^%#^&!!"
$sdeSheet.Range("A2:B2").Style = "Comma"
$sdeSheet.Range("B2").Select()

# --- New sheet: emailRecipient ---
$emailRecipientSheet = $wb.Worksheets.Add()
$emailRecipientSheet.Name = "emailRecipient"
$emailRecipientSheet.Range("A1").Value = "emailAddresses"
$emailRecipientSheet.Range("A2").Value = "AZLsaksjd@gmail.com"
$emailRecipientSheet.Range("A3").Value = "sdajdadkjhs@usgs.gov"
$emailRecipientSheet.Range("B5").Select()

# --- Fix up the "email" sheet selection (no more single-cell A3 selection) ---
$emailSheet = $wb.Worksheets.Item("email")
$emailSheet.Range("A1:B2").Select()

# --- Make emailRecipient the active tab ---
$emailRecipientSheet.Activate()
